$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6823908
$ws.Range("I98").Value = 7693746
$ws.Range("J98").Value = 2300750.2
$ws.Range("K98").Value = 7693746
$ws.Range("L98").Value = 2300750.2
$ws.Range("M98").Value = -7692248
$ws.Range("N98").Value = -2303746.2
$ws.Range("H107").Value = 820.7586
$ws.Range("I107").Value = 840.8148
$ws.Range("K107").Value = 840.8148
$ws.Range("M107").Value = 1079.1852
$ws.Range("H121").Value = 1018.125
$ws.Range("J121").Value = 1165.8334
$ws.Range("L121").Value = 3497.5002
$ws.Range("N121").Value = -6991.5002
$ws.Range("H122").Value = 6823908
$ws.Range("I122").Value = 7693746
$ws.Range("J122").Value = 2300750.2
$ws.Range("K122").Value = 23081238
$ws.Range("L122").Value = 6902250.600000001
$ws.Range("M122").Value = -23078788
$ws.Range("N122").Value = -6907150.600000001
$ws.Range("H132").Value = 1158669.5
$ws.Range("I132").Value = 1071.012
$ws.Range("J132").Value = 9261859
$ws.Range("K132").Value = 3213.036
$ws.Range("L132").Value = 27785577
$ws.Range("M132").Value = -683.0360000000001
$ws.Range("N132").Value = -27790637
$ws.Range("H138").Value = 2555.25
$ws.Range("I138").Value = 1814.4565
$ws.Range("J138").Value = 3691.1333
$ws.Range("K138").Value = 5443.3695
$ws.Range("L138").Value = 11073.3999
$ws.Range("M138").Value = -303.3694999999998
$ws.Range("N138").Value = -21353.3999
$ws.Range("H141").Value = 1341.75
$ws.Range("I141").Value = 776.18866
$ws.Range("J141").Value = 11333.333
$ws.Range("K141").Value = 2328.56598
$ws.Range("L141").Value = 33999.999
$ws.Range("M141").Value = 2851.43402
$ws.Range("N141").Value = -44359.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 477407.62
$ws.Range("I45").Value = 589408.5600000001
$ws.Range("J45").Value = 1403.5
$ws.Range("K45").Value = 589408.5600000001
$ws.Range("L45").Value = 1403.5
$ws.Range("M45").Value = -589031.5600000001
$ws.Range("N45").Value = -2157.5
$ws.Range("H74").Value = 35089108
$ws.Range("I74").Value = 37037720
$ws.Range("J74").Value = 30306146
$ws.Range("K74").Value = 37037720
$ws.Range("L74").Value = 30306146
$ws.Range("M74").Value = -37036846
$ws.Range("N74").Value = -30307894
$ws.Range("H77").Value = 35089108
$ws.Range("I77").Value = 37037720
$ws.Range("J77").Value = 30306146
$ws.Range("K77").Value = 185188600
$ws.Range("L77").Value = 151530730
$ws.Range("M77").Value = -185184232
$ws.Range("N77").Value = -151539466

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1478.0476
$ws.Range("I16").Value = 1385.2667
$ws.Range("J16").Value = 1710
$ws.Range("K16").Value = 1385.2667
$ws.Range("L16").Value = 1710
$ws.Range("M16").Value = -1098.2667
$ws.Range("N16").Value = -2284
$ws.Range("H58").Value = 2674834.5
$ws.Range("I58").Value = 966.3
$ws.Range("J58").Value = 6494646.5
$ws.Range("K58").Value = 966.3
$ws.Range("L58").Value = 6494646.5
$ws.Range("M58").Value = -763.3
$ws.Range("N58").Value = -6495052.5
$ws.Range("H107").Value = 558.8222
$ws.Range("I107").Value = 209.8
$ws.Range("J107").Value = 733.3333
$ws.Range("K107").Value = 209.8
$ws.Range("L107").Value = 733.3333
$ws.Range("M107").Value = 1710.2
$ws.Range("N107").Value = -4573.3333
$ws.Range("H113").Value = 1478.0476
$ws.Range("I113").Value = 1385.2667
$ws.Range("J113").Value = 1710
$ws.Range("K113").Value = 1385.2667
$ws.Range("L113").Value = 1710
$ws.Range("M113").Value = 784.7333000000001
$ws.Range("N113").Value = -6050
$ws.Range("H134").Value = 1668188.9
$ws.Range("I134").Value = 1572.9524
$ws.Range("J134").Value = 13334500
$ws.Range("K134").Value = 4718.857199999999
$ws.Range("L134").Value = 40003500
$ws.Range("M134").Value = -2183.857199999999
$ws.Range("N134").Value = -40008570
$ws.Range("H136").Value = 2674834.5
$ws.Range("I136").Value = 966.3
$ws.Range("J136").Value = 6494646.5
$ws.Range("K136").Value = 2898.9
$ws.Range("L136").Value = 19483939.5
$ws.Range("M136").Value = -348.8999999999996
$ws.Range("N136").Value = -19489039.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2661.75
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2661.75
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7985.25
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -8153.25
$ws.Range("H112").Value = 23711.4
$ws.Range("J112").Value = 2265
$ws.Range("L112").Value = 6795
$ws.Range("N112").Value = -9011
$ws.Range("H113").Value = 3524.5208
$ws.Range("I113").Value = 1328.5714
$ws.Range("J113").Value = 4428.7354
$ws.Range("K113").Value = 3985.7142
$ws.Range("L113").Value = 13286.2062
$ws.Range("M113").Value = -1815.7142
$ws.Range("N113").Value = -17626.2062
$ws.Range("H122").Value = 1033.6
$ws.Range("I122").Value = 264.46667
$ws.Range("J122").Value = 1610.45
$ws.Range("K122").Value = 2380.20003
$ws.Range("L122").Value = 14494.05
$ws.Range("M122").Value = 69.79997000000003
$ws.Range("N122").Value = -19394.05
$ws.Range("H131").Value = 8475406
$ws.Range("I131").Value = 41667100
$ws.Range("J131").Value = 931.34045
$ws.Range("K131").Value = 125001300
$ws.Range("L131").Value = 2794.02135
$ws.Range("M131").Value = -124996260
$ws.Range("N131").Value = -12874.02135
$ws.Range("H132").Value = 2537.12
$ws.Range("I132").Value = 2395
$ws.Range("K132").Value = 21555
$ws.Range("M132").Value = -19025

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5496557.5
$ws.Range("I132").Value = 7520078
$ws.Range("J132").Value = 4144
$ws.Range("K132").Value = 22560234
$ws.Range("L132").Value = 12432
$ws.Range("M132").Value = -22557704
$ws.Range("N132").Value = -17492
$ws.Range("H136").Value = 3585530.5
$ws.Range("I136").Value = 3969611.8
$ws.Range("J136").Value = 771
$ws.Range("K136").Value = 11908835.4
$ws.Range("L136").Value = 2313
$ws.Range("M136").Value = -11906285.4
$ws.Range("N136").Value = -7413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 789805.75
$ws.Range("J132").Value = 4116656
$ws.Range("L132").Value = 12349968
$ws.Range("N132").Value = -12355028
$ws.Range("H136").Value = 1072.8309
$ws.Range("I136").Value = 525.9459000000001
$ws.Range("J136").Value = 1667.9706
$ws.Range("K136").Value = 1577.8377
$ws.Range("L136").Value = 5003.9118
$ws.Range("M136").Value = 972.1623
$ws.Range("N136").Value = -10103.9118
